$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 745, shifting existing rows 745-786 down to 746-787.
$ws.Rows.Item(745).Insert()

# Populate the newly inserted row with its values.
# Column A holds a date-looking string that must stay plain text (as in the
# rest of the column) rather than being auto-converted into a date serial
# number by Excel. Temporarily force a text format while assigning the
# value, then clear the formatting again so the cell ends up with no
# explicit style, matching the rest of the sheet.
$dateCell = $ws.Cells.Item(745, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/31"
$dateCell.ClearFormats()

$ws.Cells.Item(745, 2).Value = "土"
$ws.Cells.Item(745, 3).Value = 14
$ws.Cells.Item(745, 4).Value = 201
